# Generate Report for Handoff
#
# Refreshes the localization-status report for the 4 files that just went
# through a handoff (2b2ddfab, 7357d7fb, 876cc178, a194e192):
#   - Priority bumped from "low" to "ht"
#   - Handoff timestamps stamped with the new generation time
#
# Overview!G4:G7 gets the new "Latest HO Xliff Generate Date"; each locale
# sheet (zh-cn, de-de) gets its own "Latest Handoff Datetime" in column H,
# plus the Priority update in column E — for data rows 4-7 (the header is
# row 1, rows 2-3 are unaffected "Handed back" files).

$wb = $excel.ActiveWorkbook

$rows = 4..7

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-22 12:30:54"
}

# --- zh-cn sheet: column E = Priority, column H = Latest Handoff Datetime ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-22 12:30:49"
}

# --- de-de sheet: column E = Priority, column H = Latest Handoff Datetime ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-22 12:30:54"
}
